$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.3
$ws.Range("G2").Value = 3.4
$ws.Range("H2").Value = 2.26
$ws.Range("I2").Value = 2.28
$ws.Range("J2").Value = 3.85
$ws.Range("K2").Value = 3.9
$ws.Range("L2").Value = 1.35
$ws.Range("R2").Value = 1.51
$ws.Range("S2").Value = 2.9
$ws.Range("V2").Value = 1.78
$ws.Range("W2").Value = 1.42
$ws.Range("Z2").Value = 15.5
$ws.Range("AI2").Value = 32
$ws.Range("AN2").Value = 26
$ws.Range("N3").Value = 3.65
$ws.Range("Q3").Value = 1.86
$ws.Range("S3").Value = 3.15
$ws.Range("U3").Value = 2.14
$ws.Range("X3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("F4").Value = 2.64
$ws.Range("G4").Value = 3.05
$ws.Range("H4").Value = 2.52
$ws.Range("I4").Value = 2.9
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 3.85
$ws.Range("L4").Value = 1.33
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 3.6
$ws.Range("O4").Value = 1.16
$ws.Range("Q4").Value = 1.88
$ws.Range("R4").Value = 1.35
$ws.Range("S4").Value = 3.25
$ws.Range("T4").Value = 1.7
$ws.Range("U4").Value = 2.14
$ws.Range("V4").Value = 1.52
$ws.Range("W4").Value = 1.49
$ws.Range("Y4").Value = 14
$ws.Range("AC4").Value = 9.6
$ws.Range("G5").Value = 1.36
$ws.Range("H5").Value = 15
$ws.Range("I5").Value = 18
$ws.Range("J5").Value = 4.9
$ws.Range("K5").Value = 5.2
$ws.Range("L5").Value = 1.42
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 3.15
$ws.Range("O5").Value = 1.37
$ws.Range("P5").Value = 1.74
$ws.Range("Q5").Value = 2.1
$ws.Range("R5").Value = 1.27
$ws.Range("S5").Value = 3.95
$ws.Range("T5").Value = 2.74
$ws.Range("U5").Value = 1.56
$ws.Range("X5").Value = 15.5
$ws.Range("AB5").Value = 6.2
$ws.Range("AC5").Value = 13.5
$ws.Range("AG5").Value = 12.5
$ws.Range("AH5").Value = 60
$ws.Range("AL5").Value = 85
$ws.Range("AN5").Value = 8.4
$ws.Range("F6").Value = 1.95
$ws.Range("G6").Value = 2.02
$ws.Range("H6").Value = 4.1
$ws.Range("I6").Value = 4.7
$ws.Range("K6").Value = 3.7
$ws.Range("Q6").Value = 1.8
$ws.Range("T6").Value = 1.71
$ws.Range("V6").Value = 1.27
$ws.Range("W6").Value = 1.98
$ws.Range("Z6").Value = 38
$ws.Range("G7").Value = 1.47
$ws.Range("W7").Value = 3.1
$ws.Range("X7").Value = 80
$ws.Range("G8").Value = 1.71
$ws.Range("H8").Value = 7.2
$ws.Range("P8").Value = 1.49
$ws.Range("W8").Value = 2.4
$ws.Range("F9").Value = 1.81
$ws.Range("G9").Value = 1.89
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = 5.8
$ws.Range("K9").Value = 3.9
$ws.Range("L9").Value = 1.38
$ws.Range("N9").Value = 3.25
$ws.Range("P9").Value = 1.76
$ws.Range("V9").Value = 1.22
$ws.Range("W9").Value = 2.12
$ws.Range("X9").Value = 24
$ws.Range("AD9").Value = 980
$ws.Range("AF9").Value = 24
$ws.Range("AK9").Value = 1000
$ws.Range("AN9").Value = 85
$ws.Range("G10").Value = 1.79
$ws.Range("H10").Value = 5.6
$ws.Range("I10").Value = 6.4
$ws.Range("J10").Value = 3.5
$ws.Range("Q10").Value = 1.84
$ws.Range("R10").Value = 1.3
$ws.Range("S10").Value = 3.6
$ws.Range("V10").Value = 1.18
$ws.Range("W10").Value = 2.26
$ws.Range("AC10").Value = 14
$ws.Range("AH10").Value = 65
$ws.Range("AN10").Value = 85
$ws.Range("G11").Value = 2.34
$ws.Range("I11").Value = 3.75
$ws.Range("Q11").Value = 2
$ws.Range("W11").Value = 1.74
$ws.Range("F12").Value = 2.38
$ws.Range("G12").Value = 2.42
$ws.Range("H12").Value = 3.35
$ws.Range("I12").Value = 3.5
$ws.Range("P12").Value = 1.74
$ws.Range("Q12").Value = 2.18
$ws.Range("R12").Value = 1.27
$ws.Range("T12").Value = 1.88
$ws.Range("W12").Value = 1.7
$ws.Range("Y12").Value = 13
$ws.Range("AA12").Value = 170
$ws.Range("AC12").Value = 7.6
$ws.Range("AG12").Value = 14
$ws.Range("AI12").Value = 65
$ws.Range("AL12").Value = 50
$ws.Range("AM12").Value = 580
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 2.1
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 4.6
$ws.Range("J13").Value = 3.4
$ws.Range("K13").Value = 3.65
$ws.Range("O13").Value = 1.38
$ws.Range("P13").Value = 1.76
$ws.Range("Q13").Value = 2.1
$ws.Range("V13").Value = 1.28
$ws.Range("W13").Value = 1.9
$ws.Range("X13").Value = 23
$ws.Range("AC13").Value = 8.199999999999999
$ws.Range("AL13").Value = 1000
